# Update "想去人数" (interested-count) values across sheets, per upstream data refresh.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 617
$wsExhibition.Range("F3").Value = 283
$wsExhibition.Range("F5").Value = 763
$wsExhibition.Range("F6").Value = 419
$wsExhibition.Range("F10").Value = 260
$wsExhibition.Range("F11").Value = 7005
$wsExhibition.Range("F14").Value = 43
$wsExhibition.Range("F16").Value = 381
$wsExhibition.Range("F22").Value = 191
$wsExhibition.Range("F23").Value = 107
$wsExhibition.Range("F27").Value = 34
$wsExhibition.Range("F28").Value = 1991
$wsExhibition.Range("F29").Value = 547
$wsExhibition.Range("F32").Value = 3

# Sheet "本地生活" (rId3 / sheet3.xml)
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 322

# Sheet "全部类型" (rId4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 322
$wsAll.Range("F3").Value = 617
$wsAll.Range("F4").Value = 283
$wsAll.Range("F6").Value = 763
$wsAll.Range("F8").Value = 419
$wsAll.Range("F12").Value = 260
$wsAll.Range("F13").Value = 7005
$wsAll.Range("F17").Value = 43
$wsAll.Range("F19").Value = 381
$wsAll.Range("F32").Value = 191
$wsAll.Range("F33").Value = 107
$wsAll.Range("F37").Value = 34
$wsAll.Range("F38").Value = 1991
$wsAll.Range("F39").Value = 547
$wsAll.Range("F42").Value = 3
